$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 42636.606678240743

$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 9991.9500000000007
$ws.Range("D5").Value = 9974
$ws.Range("E5").Value = 19.29
$ws.Range("F5").Value = 19.22

$ws.Range("G3").Copy($ws.Range("G5"))
$ws.Range("G5").Value = $true

$ws.Range("H5").Value = -0.36
$ws.Range("I5").Value = $false
